# "Update codes and figures"
#
# The source table (Sheet1, A1:E5) lists model/approach comparisons.
# Row 4, column A used to reference "Abeyshu et al. 2022" — update the
# citation text to the corrected short form "Abeshu".
#
# Also update the worksheet's last active selection to cell E8 (a cell
# just below the table), matching the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the citation/model-name cell.
$ws.Range("A4").Value = "Abeshu"

# Move/save the selection to E8 (single cell selected, not a range).
$ws.Range("E8").Select()
